# Applies the cryptos.xlsx data refresh described in the commit
# "Updated cryptos list on Tue Apr  9 04:20:14 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.117.45"
$ws.Range("E2").Value = "  +2.52%  "
$ws.Range("D3").Value = "3.682.90"
$ws.Range("E3").Value = "  +7.68%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'580.92"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").Value = "'177.38"
$ws.Range("E6").Value = "  +0.60%  "
$ws.Range("D7").Value = "3.669.41"
$ws.Range("E7").Value = "  +7.48%  "
$ws.Range("D8").Value = "'0.613"
$ws.Range("E8").Value = "  +3.72%  "
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("D11").Value = "'6.57"
$ws.Range("E11").Value = "  +21.48%  "
$ws.Range("D12").Value = "'0.606"
$ws.Range("E12").Value = "  +4.40%  "
$ws.Range("D13").Value = "'48.84"
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("D14").Value = "'0.0000285"
$ws.Range("E14").Value = "  +2.04%  "
$ws.Range("D15").Value = "4.279.79"
$ws.Range("E15").Value = "  +7.82%  "
$ws.Range("D16").Value = "'677.50"
$ws.Range("E16").Value = "  -1.87%  "
$ws.Range("D17").Value = "'8.96"
$ws.Range("E17").Value = "  +4.19%  "
$ws.Range("D18").Value = "3.691.96"
$ws.Range("E18").Value = "  +7.84%  "
$ws.Range("D19").Value = "71.268.36"
$ws.Range("E19").Value = "  +2.61%  "
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("D21").Value = "'17.88"
$ws.Range("E21").Value = "  +1.56%  "
$ws.Range("E22").Value = "  +1.52%  "
$ws.Range("D23").Value = "'0.938"
$ws.Range("E23").Value = "  +5.12%  "
$ws.Range("D24").Value = "'17.31"
$ws.Range("E24").Value = "  +2.74%  "
$ws.Range("D25").Value = "'101.94"
$ws.Range("E25").Value = "  +1.27%  "
$ws.Range("E26").Value = "  +1.97%  "
$ws.Range("D27").Value = "'2.81"
$ws.Range("E27").Value = "  +5.71%  "
$ws.Range("D28").Value = "'10.22"
$ws.Range("E28").Value = "  +7.29%  "
$ws.Range("D29").Value = "'34.93"
$ws.Range("E29").Value = "  +4.59%  "
$ws.Range("E30").Value = "  +4.93%  "
$ws.Range("D31").Value = "'7.42"
$ws.Range("E31").Value = "  +5.92%  "
$ws.Range("E32").Value = "  +10.80%  "
$ws.Range("D33").Value = "'586.24"
$ws.Range("E33").Value = "  +1.46%  "
$ws.Range("D34").Value = "'11.16"
$ws.Range("E34").Value = "  +1.60%  "
$ws.Range("E35").Value = "  +5.01%  "
$ws.Range("D36").Value = "'58.85"
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "3.672.91"
$ws.Range("E38").Value = "  +2.76%  "
$ws.Range("E39").Value = "  +4.87%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "'35.30"
$ws.Range("E40").Value = "  +1.55%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0762"
$ws.Range("E41").Value = "  +4.75%  "
$ws.Range("D42").Value = "'3.41"
$ws.Range("E42").Value = "  +5.41%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0454"
$ws.Range("E43").Value = "  +9.24%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "'2.75"
$ws.Range("E44").Value = "  +3.54%  "
$ws.Range("E45").Value = "  +5.04%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").Value = "'2.87"
$ws.Range("E46").Value = "  +8.81%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'3.36"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("E48").Value = "  +3.37%  "
$ws.Range("D49").Value = "'1.42"
$ws.Range("E49").Value = "  -2.01%  "
$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").Value = "'136.18"
$ws.Range("E51").Value = "  +2.83%  "
